# "Refactor utility rendering logic to handle errors and improve readability"
#
# The "Cab1 Generale1" utility row (row 2) is being dropped from the
# rendered table. Deleting it shifts the remaining utility rows up, so what
# was row 3 ("Cab2 Generale1") becomes row 2, and what was row 4
# ("Cab3 Generale1") becomes row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "Cab1 Generale1" - everything below shifts up one row.
$ws.Rows("2:2").Delete()

# The surviving data rows should render with the same look as the header
# row (right-aligned number style for Cabinet/Nodo, general style for
# Utenza) instead of the old per-row look-and-feel, and match the header's
# row height.
$ws.Range("A1:B1").Copy()
$ws.Range("A2:B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows("2:3").RowHeight = 19.5
